# Update workbook per commit diff:
#  - Row 2: fix reporter/observer name typo and project year text.
#  - Insert a brand-new observation record as a new row 4 ("Robust tickgnagare" /
#    Dorcatoma robusta), which pushes the two existing "Åkerväddsantennmal" rows
#    apart: the record that used to be row 4 becomes row 3, and the record that
#    used to be row 3 becomes row 5 (only their Id/Ost/Nord/Biotop-beskrivning
#    cells actually differ from one another, everything else is identical).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: name + project text corrections -------------------------------
$ws.Range("AW2").Value = "Nils Otto Nilsson"
$ws.Range("AX2").Value = "Nils Otto Nilsson"
$ws.Range("AY2").Value = "Åsums fure 2013"

# --- Make room for the new record at row 4 ---------------------------------
# This shifts the current row 4 (A4=112156959 ...) down to row 5, while the
# current row 3 (A3=112156964 ...) stays put at row 3 for now.
$ws.Rows.Item(4).Insert()

# --- Row 3 should now hold what used to be row 4's data; row 5 already does
#     hold what used to be row 3's data (courtesy of the Insert shift above).
#     Swap the 4 cells that differ between those two records.
$ws.Range("A3").Value = 112156959
$ws.Range("Q3").Value = 445824.5356392039
$ws.Range("R3").Value = 6205211.776568725
$ws.Range("AI3").Value = "i tallskogsbryn"

$ws.Range("A5").Value = 112156964
$ws.Range("Q5").Value = 445828.4356342637
$ws.Range("R5").Value = 6205165.305277914
$ws.Range("AI5").Value = "i gles tallskog"

# --- Row 4: brand-new observation record -----------------------------------
# Force text formatting first on cells whose values look numeric/date-like so
# Excel doesn't silently coerce them (e.g. "1" -> 1, "2013-05-16" -> a date
# serial number), matching the source data which stores them as plain text.
foreach ($addr in @("I4", "Y4", "Z4", "AA4", "AB4")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A4").Value = 112141528
$ws.Range("B4").Value = 4755
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 100857
$ws.Range("F4").Value = "Robust tickgnagare"
$ws.Range("G4").Value = "Dorcatoma robusta"
$ws.Range("H4").Value = "Strand, 1938"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "ex."
$ws.Range("K4").Value = "imago/adult"
$ws.Range("N4").Value = "fönsterfälla"
$ws.Range("P4").Value = "Åsums fure, delomr 19, 580 m NO om mc-banans ledningstorn, Sk"
$ws.Range("Q4").Value = 445824.022709821
$ws.Range("R4").Value = 6205170.953796315
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Skåne"
$ws.Range("U4").Value = "Kristianstad"
$ws.Range("V4").Value = "Skåne"
$ws.Range("W4").Value = "Kristianstad"
$ws.Range("Y4").Value = "2013-05-16"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AA4").Value = "2013-05-24"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AI4").Value = "i gles tallskog"
$ws.Range("AO4").Value = "på nydöd tall, delvis barklös"
$ws.Range("AQ4").Value = "Nils Otto Nilsson"
$ws.Range("AR4").Value = "NON 04741"
$ws.Range("AW4").Value = "Nils Otto Nilsson"
$ws.Range("AX4").Value = "Nils Otto Nilsson"
$ws.Range("AY4").Value = "Åsums fure 2013"
